# Added 3 songs, fixed visibility
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the 3 new cover rows (151-153) ---
$newRows = @(
    @(149, "Maroon 5 (Pentatonix", "Moves Like Jagger"),
    @(150, "Dua Lipa (Pentatonix)", "Break My Heart"),
    @(151, "Demi Levato (Pentatonix)", "Sorry Not Sorry")
)

$startRow = 151
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# --- Fix visibility: scroll the view and reselect the active cell ---
$excel.ActiveWindow.ScrollRow = 133
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A156").Select() | Out-Null
